$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-stamp the previous entry's timestamp (source data refresh rounds it to
# millisecond precision).
$ws.Range("A4").Value = 45804.43708293982

# A new price scrape was recorded for CREATINA MONOHIDRATO EN POLVO (1Kg, 12,88€),
# appended as row 5 right after the existing row 4 entry.
$ws.Range("A5").Value = 45804.44035107633
$ws.Range("A5").NumberFormat = $ws.Range("A4").NumberFormat
$ws.Range("B5").Value = "CREATINA MONOHIDRATO EN POLVO"
$ws.Range("C5").Value = "1Kg"
$ws.Range("D5").Value = "12,88€"
